# Apply the table style change on slide 6's table (Shape index 2).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$tableShape = $s.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{C733B1B8-024A-42AF-86EA-4B3A30F4B281}")
